$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1): copy the formatting from H1 (existing bold/border/
# center-top style) onto the new I1/J1 header cells, then set their text.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-8 for the new I/J columns.
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 4

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9

$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 6

$ws.Range("I6").Value = 6
$ws.Range("J6").Value = 6

$ws.Range("I7").Value = 4
$ws.Range("J7").Value = 4

$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 5
